# Remove the "[PUMP:TBD:1]" paragraph and the "BOLUS:SRS:2" (ListBullet)
# paragraph from the document, leaving the title and the blank paragraph
# that precede them, followed directly by the section properties.

$d = $word.ActiveDocument

# Locate the two paragraphs to remove by their text content so the
# script is resilient to their exact index in the document.
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.Trim()
    if ($t -eq "[PUMP:TBD:1]") {
        $startPara = $p
    }
    if ($t -eq "BOLUS:SRS:2") {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}

$d.Save()
